$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Column I = "priority" (moderate / high / low)
# Column M = "future_sampling_recommended" (yes)

$ws.Range("I4").Value = "moderate"
$ws.Range("M4").Value = "yes"

$ws.Range("I5").Value = "moderate"
$ws.Range("M5").Value = "yes"

$ws.Range("I6").Value = "high"

$ws.Range("I9").Value = "high"
$ws.Range("M9").Value = "yes"

$ws.Range("I10").Value = "high"
$ws.Range("M10").Value = "yes"

$ws.Range("I13").Value = "high"
$ws.Range("M13").Value = "yes"

$ws.Range("I15").Value = "moderate"

$ws.Range("I19").Value = "high"

$ws.Range("I20").Value = "low"

$ws.Range("I21").Value = "low"

$ws.Range("I23").Value = "low"

$ws.Range("I26").Value = "high"
$ws.Range("M26").Value = "yes"

$ws.Range("I27").Value = "moderate"
$ws.Range("M27").Value = "yes"

$ws.Range("I30").Value = "moderate"

$ws.Range("I32").Value = "moderate"

$ws.Range("I33").Value = "high"

$ws.Range("I36").Value = "high"

$ws.Range("I38").Value = "high"

$ws.Range("I41").Value = "high"

# Restore the selected cell / scroll position noted in the workbook view
$ws.Range("M41").Select()
